$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — match the formatting of the
# existing header row (copy H1's format, which carries style index 1:
# bold font, border, centered/top alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF), rows 2-8 — plain numeric values,
# no special style (matches existing unstyled data columns).
$data = @{
    2 = @(7, 7)
    3 = @(5, 6)
    4 = @(3, 4)
    5 = @(8, 8)
    6 = @(1, 3)
    7 = @(7, 8)
    8 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $pair = $data[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
